# CHARMS - Updated NV scenario code to reflect latest changes
#
# This script reproduces, via the Excel COM object model, the edits made to
# RASScenario_Consent_Adult.xlsx:
#   1. A new Question/Answer row is inserted at row 11 of the
#      "screenerScenarioAdult" sheet ("Do you currently live in the United
#      States?" / "Yes"), pushing all subsequent rows down by one.
#   2. The worksheet selection on that sheet moves to A11.
#   3. The hyperlinks that were attached to the (now shifted) rows are
#      re-attached to their new locations.
#   4. The scroll position on the "RASSurveyScenario1" sheet is adjusted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # screenerScenarioAdult

# --- 1. Insert the new row ------------------------------------------------
# Inserting a whole row shifts every cell (values, formulas, styles) below it
# down by one and keeps the sheet dimension / row spans in sync.
$ws.Rows.Item(11).Insert()

# New question cell (bold black "Calibri" 12pt, same look used elsewhere in
# the workbook for section-style headers).
$qCell = $ws.Range("A11")
$qCell.Value2 = "Do you currently live in the United States?"
$qCell.Font.Bold = $true
$qCell.Font.Size = 12
$qCell.Font.Color = 0

# New answer cell ("Yes", right-aligned, regular "Calibri" 11pt black).
$aCell = $ws.Range("B11")
$aCell.Value2 = "Yes"
$aCell.Font.Bold = $false
$aCell.Font.Size = 11
$aCell.Font.Color = 0
$aCell.HorizontalAlignment = -4152   # xlRight

# --- 2. Fix up the hyperlinks ----------------------------------------------
# The two mailto hyperlinks used to sit on B18/B19; after the row insert the
# underlying cell data now lives on B19/B20. Preserve their original
# formatting (the "Hyperlink" look, cell style index) by stashing a copy of
# it in an unused scratch cell before Hyperlinks.Add() re-applies its own
# default styling, then paste the formatting back and discard the scratch
# cell via a row delete (which keeps the sheet dimension unchanged).
$scratchRow = $ws.Rows.Item($ws.Rows.Count)
$scratchAddr = "B" + ($ws.UsedRange.Rows.Count + 1)
$ws.Range("B19").Copy($ws.Range($scratchAddr))

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:consent_participant@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("B20"), "mailto:consent_participant@yopmail.com")

$ws.Range($scratchAddr).Copy()
$ws.Range("B19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$scratchRowNum = [int]($scratchAddr -replace '[^0-9]', '')
$ws.Rows.Item($scratchRowNum).Delete()

# --- 3. Update the selection on this sheet ----------------------------------
$ws.Range("A11").Select()

# --- 4. Adjust the scroll position on RASSurveyScenario1 -------------------
$ws4 = $wb.Worksheets.Item(4)   # RASSurveyScenario1
$ws4.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$win.ScrollColumn = 1

# Restore the originally active sheet/tab.
$ws.Activate()
